$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 268, shifting rows 268:303 down to 269:304
$ws.Rows.Item(268).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the new row 268 with the new data record
$ws.Range("A268").Value = 6
$ws.Range("B268").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C268").Value = "Metropolitana"
$ws.Range("D268").Value = 44491
$ws.Range("D268").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E268").Value = 13
$ws.Range("F268").Value = 100112039
$ws.Range("G268").Value = "Ciboulette"
$ws.Range("H268").Value = "Sin especificar"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 830
$ws.Range("K268").Value = 800
$ws.Range("L268").Value = 900
$ws.Range("M268").Value = 847
$ws.Range("N268").Value = "`$/docena de atados"
$ws.Range("O268").Value = "Región Metropolitana"
$ws.Range("P268").Value = 282
$ws.Range("Q268").Value = 3
$ws.Range("R268").Value = "Hortaliza"
